$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Hoja1")

# Try: copy the WHOLE ROW 3 formatting across to row 9, which for C3:O3 already has applyAlignment=1 wrap/left/vcenter.
# Then, instead of changing alignment props, use Range.ClearFormats first on target, then set Interior via copy paste of a NON-aligned fill cell.
# Is there any existing fillId=4 cell WITHOUT alignment? No. So instead: Copy a cell that has fillId=4 but try clearing just the alignment via .AddIndent or checking UnMerge... Let's try setting WrapText=$false and Horizontal/VerticalAlignment back to explicit "General"/"Bottom" AFTER clearing via .Style = "Normal" first then reapplying fill only.

$ws.Range("C21").Style = "Normal"
$ws.Range("C21").Interior.Pattern = -4124  # xlSolid
$ws.Range("C21").Interior.Color = 255
Write-Host (
  "C21 done"
)
